# Week 13 logging update for the Chiefs 2021 Team Data workbook.
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# YDS sheet: append this week's per-play yardage logs (Rushing/Passing,
# Offense/Defense) to the four running logs.
# ---------------------------------------------------------------------------
$ydsWs = $wb.Worksheets.Item("YDS")

$ydsWs.Cells.Item(2, 2).Value2 = $ydsWs.Cells.Item(2, 2).Value2 + " 10 2 3 5 10 4 2 12 4 2 0 5 4 2 4 4 4 4 6 7 -4 1"
$ydsWs.Cells.Item(3, 2).Value2 = $ydsWs.Cells.Item(3, 2).Value2 + " 9 9 12 13 9 29 14 -7 12 6 14 14 38 8 4"
$ydsWs.Cells.Item(2, 3).Value2 = $ydsWs.Cells.Item(2, 3).Value2 + " 9 -1 1 14 2 9 6 7 2 14 3 0 1 11 3 5 5 -1 5 4 8 9 7 2 6 2 16 1"
$ydsWs.Cells.Item(3, 3).Value2 = $ydsWs.Cells.Item(3, 3).Value2 + " 36 10 10 6 3 2 34 5 9 5 8 11 9 13 19 12 6 19 10 7"

# ---------------------------------------------------------------------------
# OFF sheet: updated situational-down counters for Home/Road.
# ---------------------------------------------------------------------------
$offWs = $wb.Worksheets.Item("OFF")

$offWs.Cells.Item(2, 2).Value2 = 2
$offWs.Cells.Item(2, 3).Value2 = 148
$offWs.Cells.Item(2, 5).Value2 = 5
$offWs.Cells.Item(2, 6).Value2 = 35
$offWs.Cells.Item(2, 7).Value2 = 41
$offWs.Cells.Item(2, 9).Value2 = 5
$offWs.Cells.Item(2, 10).Value2 = 17
$offWs.Cells.Item(2, 12).Value2 = 286
$offWs.Cells.Item(2, 13).Value2 = 174
$offWs.Cells.Item(2, 17).Value2 = 473

$offWs.Cells.Item(3, 2).Value2 = 11
$offWs.Cells.Item(3, 3).Value2 = 163
$offWs.Cells.Item(3, 5).Value2 = 27
$offWs.Cells.Item(3, 6).Value2 = 88
$offWs.Cells.Item(3, 7).Value2 = 42
$offWs.Cells.Item(3, 8).Value2 = 20
$offWs.Cells.Item(3, 9).Value2 = 41
$offWs.Cells.Item(3, 10).Value2 = 44
$offWs.Cells.Item(3, 14).Value2 = 13

# ---------------------------------------------------------------------------
# DEF sheet: updated situational-down counters for Home/Road.
# ---------------------------------------------------------------------------
$defWs = $wb.Worksheets.Item("DEF")

$defWs.Cells.Item(2, 3).Value2 = 129
$defWs.Cells.Item(2, 4).Value2 = 8
$defWs.Cells.Item(2, 6).Value2 = 51
$defWs.Cells.Item(2, 7).Value2 = 41
$defWs.Cells.Item(2, 9).Value2 = 5
$defWs.Cells.Item(2, 10).Value2 = 24
$defWs.Cells.Item(2, 12).Value2 = 242
$defWs.Cells.Item(2, 13).Value2 = 153
$defWs.Cells.Item(2, 15).Value2 = 21
$defWs.Cells.Item(2, 16).Value2 = 13
$defWs.Cells.Item(2, 17).Value2 = 423

$defWs.Cells.Item(3, 2).Value2 = 11
$defWs.Cells.Item(3, 3).Value2 = 141
$defWs.Cells.Item(3, 5).Value2 = 22
$defWs.Cells.Item(3, 6).Value2 = 70
$defWs.Cells.Item(3, 7).Value2 = 27
$defWs.Cells.Item(3, 8).Value2 = 15
$defWs.Cells.Item(3, 9).Value2 = 43
$defWs.Cells.Item(3, 10).Value2 = 34
$defWs.Cells.Item(3, 14).Value2 = 13

# ---------------------------------------------------------------------------
# ST sheet: kicking/special-teams totals and per-week logs.
# ---------------------------------------------------------------------------
$stWs = $wb.Worksheets.Item("ST")

$stWs.Cells.Item(2, 2).Value2 = 64
$stWs.Cells.Item(2, 4).Value2 = 27
$stWs.Cells.Item(2, 6).Value2 = 234
$stWs.Cells.Item(2, 7).Value2 = 219
$stWs.Cells.Item(2, 10).Value2 = 87
$stWs.Cells.Item(2, 11).Value2 = 85
$stWs.Cells.Item(2, 12).Value2 = 42
$stWs.Cells.Item(2, 13).Value2 = 36
$stWs.Cells.Item(2, 14).Value2 = 26
$stWs.Cells.Item(2, 15).Value2 = 18

$stWs.Cells.Item(3, 2).Value2 = 40

$stWs.Cells.Item(3, 4).Value2 = $stWs.Cells.Item(3, 4).Value2 + " 54 50 37"
$stWs.Cells.Item(4, 2).Value2 = $stWs.Cells.Item(4, 2).Value2 + " 66 63"
$stWs.Cells.Item(4, 4).Value2 = $stWs.Cells.Item(4, 4).Value2 + " 14 6 0"
$stWs.Cells.Item(5, 2).Value2 = $stWs.Cells.Item(5, 2).Value2 + " 26 18"
$stWs.Cells.Item(5, 4).Value2 = "0 0"

# ---------------------------------------------------------------------------
# TURNS sheet: turnover counts for Home/Road.
# ---------------------------------------------------------------------------
$turnsWs = $wb.Worksheets.Item("TURNS")

$turnsWs.Cells.Item(2, 2).Value2 = 8
$turnsWs.Cells.Item(2, 3).Value2 = 8
$turnsWs.Cells.Item(2, 5).Value2 = 7

$turnsWs.Cells.Item(3, 4).Value2 = 7
$turnsWs.Cells.Item(3, 5).Value2 = 6

# ---------------------------------------------------------------------------
# PEN sheet: penalty counts.
# ---------------------------------------------------------------------------
$penWs = $wb.Worksheets.Item("PEN")

$penWs.Cells.Item(2, 2).Value2 = 14
$penWs.Cells.Item(2, 4).Value2 = 6

$penWs.Cells.Item(3, 2).Value2 = 18

$penWs.Cells.Item(4, 2).Value2 = 4

# ---------------------------------------------------------------------------
# Leave the ST sheet active with D6 selected, matching the author's
# last-saved UI state.
# ---------------------------------------------------------------------------
$stWs.Activate()
$stWs.Range("D6").Select()
